$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before column A. This shifts the existing
# A/B/C data (segment names / PercActivations / PercSegmentAreas),
# header texts included, one column to the right -> B/C/D.
$ws.Columns.Item(1).Insert()

# New header cell B1 ("segments") - copy the header-row formatting
# (bold + border + centered) from the neighboring header cell C1, then
# set its text.
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B1").Value = "segments"

# The numeric-index column A2:A20 takes over the header-style formatting
# that used to live on the names column (bold + border + center/top).
$ws.Range("C1").Copy()
$ws.Range("A2:A20").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# The names column (now B2:B20) becomes plain, unstyled data - matching
# how the numeric columns were unstyled before the edit.
$ws.Range("B2:B20").ClearFormats()

# Fill the new column A with the 0-based segment index.
for ($i = 0; $i -lt 19; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $i
}

Write-Output "done"
